# Commit: "feat: add 2022-Q4 data"
#
# The workbook tracks 00762 (China Unicom) holdings by China-fund quarterly
# reports. A new quarter, 2022-Q4, is added:
#   1. A new worksheet "2022-Q4" is inserted right after "总计", holding the
#      24 individual-fund rows for that quarter (pushing every existing
#      quarter sheet one slot later in the tab order - that shift is handled
#      automatically because we only ever address sheets by name).
#   2. The "总计" (summary) sheet gets a new row 2 for 2022-Q4
#      (24 holders, 5.52 billion yuan held), with all prior rows' running
#      index (column A) bumped by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet right after "总计", copying the
#    formatting (header style, borders, column layout) from the existing
#    "2022-Q3" sheet so the new sheet matches the others' look.
# ---------------------------------------------------------------------
$zong = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2022-Q3")
$template.Copy($null, $zong)

# The copy lands immediately after "总计"; it is named "2022-Q3 (2)".
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# Clear out the template's old data rows (2-9) before writing the real
# 2022-Q4 numbers (which span rows 2-25).
$newSheet.Range("A2:H9").ClearContents()

# ---------------------------------------------------------------------
# 2. Populate "2022-Q4" fund-level holdings.
# ---------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$rows = @(
    @(0,  "014887", "招商安福1年定期开放债券",              "17.78", "33.59", "6.75", "1.2002", 1),
    @(1,  "002121", "广发沪港深新起点股票A",                "26.30", "88.97", "3.62", "0.9521", 8),
    @(2,  "010430", "招商安阳债券A",                        "16.90", "20.35", "4.42", "0.7470", 1),
    @(3,  "016513", "招商安嘉债券",                          "16.27", "20.17", "4.59", "0.7468", 1),
    @(4,  "010591", "富国中国中小盘混合（QDII）美元",        "32.21", "87.21", "2.21", "0.7118", 8),
    @(5,  "100061", "富国中国中小盘混合（QDII）人民币",      "32.21", "87.21", "2.21", "0.7118", 8),
    @(6,  "012770", "光大保德信创新生活混合",                "2.89",  "88.40", "4.43", "0.1280", 5),
    @(7,  "014094", "南方誉盈一年持有混合A",                "8.21",  "28.74", "1.03", "0.0846", 9),
    @(8,  "005143", "中融沪港深大消费主题灵活配置混合C",    "0.69",  "77.26", "5.82", "0.0402", 2),
    @(9,  "007107", "太平 MSCI 香港价值增强指数A",          "1.02",  "92.16", "2.48", "0.0253", 9),
    @(10, "014697", "南方誉稳一年持有混合A",                "3.00",  "24.51", "0.84", "0.0252", 4),
    @(11, "014698", "南方誉稳一年持有混合C",                "2.83",  "24.51", "0.84", "0.0238", 4),
    @(12, "014214", "光大保德信核心资产混合A",              "0.32",  "85.79", "7.14", "0.0228", 1),
    @(13, "010024", "广发沪港深新起点股票C",                "0.59",  "88.97", "3.62", "0.0214", 8),
    @(14, "005142", "中融沪港深大消费主题灵活配置混合A",    "0.31",  "77.26", "5.82", "0.0180", 2),
    @(15, "014146", "景顺长城港股通数字经济主题混合A",      "0.60",  "89.97", "2.76", "0.0166", 8),
    @(16, "014462", "光大保德信汇佳混合A",                  "0.27",  "88.46", "5.35", "0.0144", 2),
    @(17, "005269", "华泰柏瑞港股通量化灵活配置混合",      "0.54",  "80.96", "1.89", "0.0102", 6),
    @(18, "014147", "景顺长城港股通数字经济主题混合C",      "0.29",  "89.97", "2.76", "0.0080", 8),
    @(19, "014095", "南方誉盈一年持有混合C",                "0.53",  "28.74", "1.03", "0.0055", 9),
    @(20, "010431", "招商安阳债券C",                        "0.11",  "20.35", "4.42", "0.0049", 1),
    @(21, "014215", "光大保德信核心资产混合C",              "0.03",  "85.79", "7.14", "0.0021", 1),
    @(22, "014463", "光大保德信汇佳混合C",                  "0.01",  "88.46", "5.35", "0.0005", 2),
    @(23, "007108", "太平 MSCI 香港价值增强指数C",          "0.00",  "92.16", "2.48", 0,        9)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3. Update the "总计" summary sheet: insert a new row for 2022-Q4 at the
#    top of the data (row 2), pushing the existing rows down and bumping
#    their running index in column A.
# ---------------------------------------------------------------------
$zong.Rows.Item(2).Insert()

# Re-apply the data-row formatting (the auto-inserted row borrows the
# header's bold style) by pasting the format from the row right below.
$zong.Range("A3:D3").Copy()
$zong.Range("A2:D2").PasteSpecial(-4122)

$zong.Range("A2").Value = 0
$zong.Range("B2").Value = "2022-Q4"
$zong.Range("C2").Value = 24
$zong.Range("D2").Value = 5.52

for ($i = 0; $i -le 7; $i++) {
    $zong.Cells.Item(3 + $i, 1).Value = $i + 1
}
